$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.455
$ws.Range("E4").Value = 13.421

$ws.Range("E5").Value = 13.41

$ws.Range("D6").Value = -8.181999999999999
$ws.Range("E6").Value = 12.805

$ws.Range("D7").Value = -7.703

$ws.Range("D8").Value = -7.779000000000001
$ws.Range("E8").Value = 13.486

$ws.Range("D16").Value = -7.783999999999999
$ws.Range("E16").Value = 12.908

$ws.Range("D20").Value = -8.093

$ws.Range("D21").Value = -8.1

$ws.Range("E22").Value = 13.406
